# Update cryptocurrency price/volume data per the Apr 27 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.361.75'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '3.237.60'
$ws.Range('E3').Value = '  +2.87%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '''594.91'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.34%  '
$ws.Range('D6').Value = '''141.76'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.11%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '3.230.89'
$ws.Range('E8').Value = '  +2.88%  '
$ws.Range('E9').Value = '  -1.54%  '
$ws.Range('E10').Value = '  -1.41%  '
$ws.Range('D11').Value = '''5.33'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.99%  '
$ws.Range('E12').Value = '  -0.75%  '
$ws.Range('D13').Value = '''0.0000247'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.00%  '
$ws.Range('D14').Value = '''34.41'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.82%  '
$ws.Range('D15').Value = '3.770.69'
$ws.Range('E15').Value = '  +2.98%  '
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = '3.238.87'
$ws.Range('E17').Value = '  +3.02%  '
$ws.Range('D18').Value = '63.402.28'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('D19').Value = '''6.78'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.36%  '
$ws.Range('D20').Value = '''475.05'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.94%  '
$ws.Range('E21').Value = '  -3.62%  '
$ws.Range('D22').Value = '''0.727'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.90%  '
$ws.Range('D23').Value = '''7.92'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.00%  '
$ws.Range('D24').Value = '''84.04'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -4.85%  '
$ws.Range('D25').Value = '''13.17'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.20%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').Value = '''7.55'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +7.25%  '
$ws.Range('E28').Value = '  -1.09%  '
$ws.Range('E29').Value = '  -1.37%  '
$ws.Range('D30').Value = '''2.10'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.66%  '
$ws.Range('D31').Value = '''27.49'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.98%  '
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('E33').Value = '  -3.71%  '
$ws.Range('E34').Value = '  -4.74%  '
$ws.Range('E35').Value = '  -1.44%  '
$ws.Range('D36').Value = '''5.91'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.74%  '
$ws.Range('D37').Value = '''52.64'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('D38').Value = '0.0₃0710'
$ws.Range('E38').Value = '  -5.20%  '
$ws.Range('D39').Value = '''0.0393'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.30%  '
$ws.Range('D40').Value = '''423.40'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.20%  '
$ws.Range('B41').Value = 'Cosmos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D41').Value = '''8.39'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = '''2.76'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -7.14%  '
$ws.Range('D43').Value = '2.974.07'
$ws.Range('E43').Value = '  +1.46%  '
$ws.Range('E44').Value = '  -8.62%  '
$ws.Range('D45').Value = '''0.267'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.10%  '
$ws.Range('D46').Value = '''2.16'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.76%  '
$ws.Range('D47').Value = '''2.37'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.52%  '
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('D49').Value = '''25.95'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.19%  '
$ws.Range('E50').Value = '  -0.57%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = '''120.89'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.43%  '
